# The deck has a duplicate slide: slide 14 (the picture-background version of
# "The crucial issue is not about asking questions, It is the spirit in which
# questions are asked.") is a repeat of slide 15 (the solid-fill version of the
# same quote). Delete the repeated slide (slide 14), which shifts slide 15
# (and everything after it) up by one position.
$p = $ppt.ActivePresentation
$p.Slides.Item(14).Delete()
